$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(341, 18, 1, "박근혜", 49.4, 41235),
    @(342, 18, 2, "문재인", 43.2, 41235),
    @(343, 18, 3, "박종선", 0.4, 41235),
    @(344, 18, 1, "박근혜", 45, 41234),
    @(345, 18, 2, "문재인", 46, 41234),
    @(346, 18, 3, "박종선", 0.2, 41234),
    @(347, 18, 1, "박근혜", 43.4, 41237),
    @(348, 18, 2, "문재인", 37.6, 41237),
    @(349, 18, 3, "박종선", 0.1, 41237),
    @(350, 18, 1, "박근혜", 48.9, 41247),
    @(351, 18, 2, "문재인", 42.1, 41247),
    @(352, 18, 3, "박종선", 0.1, 41247),
    @(353, 18, 1, "박근혜", 48.9, 41246),
    @(354, 18, 2, "문재인", 42.8, 41246),
    @(355, 18, 3, "박종선", 0.2, 41246),
    @(356, 18, 1, "박근혜", 43.5, 41249),
    @(357, 18, 2, "문재인", 43.3, 41249),
    @(358, 18, 3, "박종선", 0.1, 41249),
    @(359, 18, 1, "박근혜", 42.9, 41255),
    @(360, 18, 2, "문재인", 43.7, 41255),
    @(361, 18, 3, "박종선", 0.2, 41255),
    @(362, 18, 1, "박근혜", 44.3, 41232),
    @(363, 18, 2, "문재인", 43.3, 41232),
    @(364, 18, 3, "박종선", 0.2, 41232),
    @(365, 18, 1, "박근혜", 48.8, 41233),
    @(366, 18, 2, "문재인", 43.3, 41233),
    @(367, 18, 3, "박종선", 0.1, 41233),
    @(368, 18, 1, "박근혜", 44.6, 41231),
    @(369, 18, 2, "문재인", 39.9, 41231),
    @(370, 18, 3, "박종선", 0, 41231)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $prevRow = $rowNum - 1
    $ws.Cells.Item($prevRow, 5).Copy()
    $ws.Cells.Item($rowNum, 5).PasteSpecial(-4122)
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
}

$ws.Range("E345").Select()

Write-Host "done"
